$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("items")

$ws.Range("C2").Value = 3
$ws.Range("D2").Value = 4
$ws.Range("F2").Value = 0

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 4

$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 4

$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 4

$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 4
$ws.Range("G6").Value = 3

$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 3

$ws.Range("C8").Value = 3
$ws.Range("D8").Value = 3

$ws.Range("C9").Value = 3
$ws.Range("D9").Value = 3

$ws.Range("C10").Value = 3
$ws.Range("D10").Value = 3

$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 3

$ws.Range("C12").Value = 3
$ws.Range("E12").Value = 3

$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 3

$ws.Range("C16").Value = 10
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 3

$ws.Range("D17").Value = 5

$ws.Range("D12").Select()
